# 20250121: Cambio en cantidad = 0
# Populate "Sheet1" with the MBI operations report (header + 7 trade rows)
# and apply the header / date-column formatting used by the source report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1)
# ---------------------------------------------------------------------------
$headers = @(
    "nombre_fondo",
    "fecha_pago",
    "fecha_ingreso",
    "precio",
    "cantidad",
    "monto",
    "comision",
    "nemotecnico",
    "compra/venta/vencimiento",
    "tipo_operacion",
    "precio_factura",
    "cantidad_factura"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2. Data rows (rows 2-8)
# ---------------------------------------------------------------------------
$rows = @(
    @("FONDO DE INVERSION NEVASA AHORRO", 45678, 45678, 0.57, 176126896,       174076278, 0, "SMT_24032025_21012025_0.57_CFIMBIDT-B",  "COMPRA", "SIMULTANEA", 15909,  10942),
    @("FONDO DE INVERSION NEVASA AHORRO", 45678, 45678, 0.57, 183128353,       181918595, 0, "SMT_25022025_21012025_0.57_BCI",         "COMPRA", "SIMULTANEA", 29177,  6235),
    @("FONDO DE INVERSION NEVASA AHORRO", 45678, 45678, 0.57, 334843162,       330944625, 0, "SMT_24032025_21012025_0.57_CFIMBIDA-A",  "COMPRA", "SIMULTANEA", 17625,  18777),
    @("FONDO DE INVERSION NEVASA AHORRO", 45678, 45678, 0.57, 72769944.848,    47499964,  0, "SMT_18022025_21012025_0.57_ENELAM",      "COMPRA", "SIMULTANEA", 91.38,  519807),
    @("FONDO DE INVERSION NEVASA AHORRO", 45678, 45678, 0.57, 216637340.952,   141408186, 0, "SMT_18022025_21012025_0.57_MALLPLAZA",   "COMPRA", "SIMULTANEA", 1654.4, 85474),
    @("FONDO DE INVERSION NEVASA AHORRO", 45678, 45678, 0.57, 72769906.54800001, 47499939,  0, "SMT_18022025_21012025_0.57_CHILE",    "COMPRA", "SIMULTANEA", 123.09, 385896),
    @("FONDO DE INVERSION NEVASA AHORRO", 45678, 45678, 0.57, 124689480,       81390000,  0, "SMT_18022025_21012025_0.57_RIPLEY",      "COMPRA", "SIMULTANEA", 271.3,  300000)
)

$rowIndex = 2
foreach ($row in $rows) {
    for ($col = 0; $col -lt $row.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col + 1).Value = $row[$col]
    }
    $rowIndex++
}

# ---------------------------------------------------------------------------
# 3. Formatting
# ---------------------------------------------------------------------------
# 3a. Header style: bold font, thin box border, centered / top-aligned.
#     Build the combined style once on a scratch cell, then paste the
#     resulting *format only* onto the header row so the style table ends
#     up with a single extra cellXf (instead of one per incremental change).
$scratchHeader = $ws.Range("Z1")
$scratchHeader.Value = "x"
$scratchHeader.Font.Bold = $true
$scratchHeader.Borders.LineStyle = 1
$scratchHeader.HorizontalAlignment = -4108
$scratchHeader.VerticalAlignment = -4160

$headerRange = $ws.Range("A1:L1")
$scratchHeader.Copy()
$headerRange.PasteSpecial(-4122)
$scratchHeader.Clear()

# 3b. Date columns (fecha_pago / fecha_ingreso): same trick, built once on a
#     scratch cell so only a single new cellXf (numFmtId only) is produced.
$scratchDate = $ws.Range("Z2")
$scratchDate.Value = 45678
$scratchDate.NumberFormat = "yyyy-mm-dd"
$scratchDate.NumberFormat = "YYYY-MM-DD"

$dateRange = $ws.Range("B2:C8")
$scratchDate.Copy()
$dateRange.PasteSpecial(-4122)
$scratchDate.Clear()

Write-Host "Populated Sheet1 with 20250121 MBI operations report"
